$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date updated
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher is now populated, and the old duplicated "Contact" row is replaced
# by a "Jurisdiction" row. This removes one row overall (21 -> 20 rows), so
# delete row 11 (the second, duplicate "Contact" row) first, then update
# row 9 (Publisher) and row 10 (was "Contact", becomes "Jurisdiction").
$meta.Rows.Item(11).Delete()

$meta.Range("A9").Value = "Publisher"
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension element's Short/Definition now describe the actual extension
$elements.Range("K2").Value = "Identified Pronouns"
$elements.Range("L2").Value = "Gender, typically described in terms of masculinity and femininity, is a social construction that varies across different cultures and over time.  This value is often used for identity purposes and should be collected directly from the patient."
